$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1226.3226
$ws.Range("J19").Value = 1678.75
$ws.Range("L19").Value = 1678.75
$ws.Range("N19").Value = -2028.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2324.8333
$ws.Range("I40").Value = 2336.75
$ws.Range("J40").Value = 2301
$ws.Range("K40").Value = 2336.75
$ws.Range("L40").Value = 2301
$ws.Range("M40").Value = -2161.75
$ws.Range("N40").Value = -2651

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 5500
$ws.Range("I46").Value = 3250
$ws.Range("K46").Value = 9750
$ws.Range("M46").Value = -9631

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 5500
$ws.Range("I60").Value = 3250
$ws.Range("K60").Value = 9750
$ws.Range("M60").Value = -9266

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3638.6
$ws.Range("I112").Value = 5500
$ws.Range("J112").Value = 3352.2307
$ws.Range("K112").Value = 16500
$ws.Range("L112").Value = 10056.6921
$ws.Range("M112").Value = -15392
$ws.Range("N112").Value = -12272.6921

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3746.0195
$ws.Range("I138").Value = 1405.25
$ws.Range("J138").Value = 4816.086
$ws.Range("K138").Value = 4215.75
$ws.Range("L138").Value = 14448.258
$ws.Range("M138").Value = 924.25
$ws.Range("N138").Value = -24728.258

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7870.381
$ws.Range("I2").Value = 6373.9165
$ws.Range("J2").Value = 9865.666999999999
$ws.Range("K2").Value = 6373.9165
$ws.Range("L2").Value = 9865.666999999999
$ws.Range("M2").Value = -6260.9165
$ws.Range("N2").Value = -10091.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 4500
$ws.Range("I6").Value = 4500
$ws.Range("K6").Value = 4500
$ws.Range("M6").Value = -4327

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2417458.8
$ws.Range("I32").Value = 1539.5278
$ws.Range("K32").Value = 1539.5278
$ws.Range("M32").Value = -1252.5278

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5628.775
$ws.Range("I61").Value = 7058.222
$ws.Range("K61").Value = 7058.222
$ws.Range("M61").Value = -6846.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4323.3887
$ws.Range("I74").Value = 4542.923
$ws.Range("J74").Value = 3752.6
$ws.Range("K74").Value = 4542.923
$ws.Range("L74").Value = 3752.6
$ws.Range("M74").Value = -3668.923
$ws.Range("N74").Value = -5500.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4323.3887
$ws.Range("I77").Value = 4542.923
$ws.Range("J77").Value = 3752.6
$ws.Range("K77").Value = 22714.615
$ws.Range("L77").Value = 18763
$ws.Range("M77").Value = -18346.615
$ws.Range("N77").Value = -27499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 36741.668
$ws.Range("J104").Value = 36741.668
$ws.Range("L104").Value = 36741.668
$ws.Range("N104").Value = -43729.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 7870.381
$ws.Range("I116").Value = 6373.9165
$ws.Range("J116").Value = 9865.666999999999
$ws.Range("K116").Value = 6373.9165
$ws.Range("L116").Value = 9865.666999999999
$ws.Range("M116").Value = -4079.9165
$ws.Range("N116").Value = -14453.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 704098.4
$ws.Range("I132").Value = 842843.6
$ws.Range("J132").Value = 114431
$ws.Range("K132").Value = 2528530.8
$ws.Range("L132").Value = 343293
$ws.Range("M132").Value = -2526000.8
$ws.Range("N132").Value = -348353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5628.775
$ws.Range("I136").Value = 7058.222
$ws.Range("K136").Value = 21174.666
$ws.Range("M136").Value = -18624.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7870.381
$ws.Range("I3").Value = 6373.9165
$ws.Range("J3").Value = 9865.666999999999
$ws.Range("K3").Value = 6373.9165
$ws.Range("L3").Value = 9865.666999999999
$ws.Range("M3").Value = -6259.9165
$ws.Range("N3").Value = -10093.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7520.2554
$ws.Range("I99").Value = 7461.4
$ws.Range("J99").Value = 7856.5713
$ws.Range("K99").Value = 7461.4
$ws.Range("L99").Value = 7856.5713
$ws.Range("M99").Value = -5963.4
$ws.Range("N99").Value = -10852.5713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 815331.1
$ws.Range("I134").Value = 1017977.44
$ws.Range("K134").Value = 3053932.32
$ws.Range("M134").Value = -3051397.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 106000
$ws.Range("J137").Value = 106000
$ws.Range("L137").Value = 106000
$ws.Range("N137").Value = -116200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 84548.17999999999
$ws.Range("J140").Value = 84548.17999999999
$ws.Range("L140").Value = 84548.17999999999
$ws.Range("N140").Value = -94908.17999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 8763.362999999999
$ws.Range("I3").Value = 9232.833000000001
$ws.Range("J3").Value = 8200
$ws.Range("K3").Value = 9232.833000000001
$ws.Range("L3").Value = 8200
$ws.Range("M3").Value = -9119.833000000001
$ws.Range("N3").Value = -8426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3120.5
$ws.Range("I31").Value = 827.25
$ws.Range("J31").Value = 5590.154
$ws.Range("K31").Value = 827.25
$ws.Range("L31").Value = 5590.154
$ws.Range("M31").Value = -532.25
$ws.Range("N31").Value = -6180.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3120.5
$ws.Range("I34").Value = 827.25
$ws.Range("J34").Value = 5590.154
$ws.Range("K34").Value = 827.25
$ws.Range("L34").Value = 5590.154
$ws.Range("M34").Value = -625.25
$ws.Range("N34").Value = -5994.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 52641010
$ws.Range("J58").Value = 13472
$ws.Range("L58").Value = 13472
$ws.Range("N58").Value = -13878

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 52641010
$ws.Range("J136").Value = 13472
$ws.Range("L136").Value = 40416
$ws.Range("N136").Value = -45516

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 534.05554
$ws.Range("I14").Value = 534.05554
$ws.Range("K14").Value = 1602.16662
$ws.Range("M14").Value = -1429.16662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1201.6666
$ws.Range("I50").Value = 1287
$ws.Range("K50").Value = 3861
$ws.Range("M50").Value = -3380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 1201.6666
$ws.Range("I53").Value = 1287
$ws.Range("K53").Value = 3861
$ws.Range("M53").Value = -3380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 167729.5
$ws.Range("I114").Value = 786.6667
$ws.Range("J114").Value = 334672.34
$ws.Range("K114").Value = 2360.0001
$ws.Range("L114").Value = 1004017.02
$ws.Range("M114").Value = 893.9998999999998
$ws.Range("N114").Value = -1010525.02

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 34671036
$ws.Range("J131").Value = 41679430
$ws.Range("L131").Value = 125038290
$ws.Range("N131").Value = -125048370

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 24511512
$ws.Range("I140").Value = 27779306
$ws.Range("J140").Value = 3062.5
$ws.Range("K140").Value = 83337918
$ws.Range("L140").Value = 9187.5
$ws.Range("M140").Value = -83332738
$ws.Range("N140").Value = -19547.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 31000
$ws.Range("I52").Value = 30000
$ws.Range("J52").Value = 35000
$ws.Range("K52").Value = 30000
$ws.Range("L52").Value = 35000
$ws.Range("M52").Value = -29741
$ws.Range("N52").Value = -35518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 61499.5
$ws.Range("J100").Value = 61499.5
$ws.Range("L100").Value = 61499.5
$ws.Range("N100").Value = -63663.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 53600
$ws.Range("J133").Value = 64333.332
$ws.Range("L133").Value = 64333.332
$ws.Range("N133").Value = -74453.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 100
$ws.Range("K18").Value = 100
$ws.Range("M18").Value = 72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 12501282
$ws.Range("I46").Value = 977.6316
$ws.Range("K46").Value = 977.6316
$ws.Range("M46").Value = -789.6316

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4351.5386
$ws.Range("I132").Value = 3676.6667
$ws.Range("J132").Value = 5870
$ws.Range("K132").Value = 11030.0001
$ws.Range("L132").Value = 17610
$ws.Range("M132").Value = -8500.000100000001
$ws.Range("N132").Value = -22670

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 9499
$ws.Range("I3").Value = 9499
$ws.Range("K3").Value = 9499
$ws.Range("M3").Value = -9385

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7047.484
$ws.Range("I132").Value = 6397.3335
$ws.Range("J132").Value = 11436
$ws.Range("K132").Value = 19192.0005
$ws.Range("L132").Value = 34308
$ws.Range("M132").Value = -16662.0005
$ws.Range("N132").Value = -39368

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 94899
$ws.Range("J135").Value = 94899
$ws.Range("L135").Value = 94899
$ws.Range("N135").Value = -106359.18
